$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 '63.285.89'
Set-TextValue 2 5 '  -7.36%  '

Set-TextValue 3 4 '3.272.78'
Set-TextValue 3 5 '  -9.09%  '

Set-TextValue 4 4 '0.999'
Set-TextValue 4 5 '  -0.35%  '

Set-TextValue 5 4 '178.03'
Set-TextValue 5 5 '  -13.91%  '

Set-TextValue 6 4 '511.34'
Set-TextValue 6 5 '  -10.17%  '

Set-TextValue 7 4 '0.589'
Set-TextValue 7 5 '  -3.80%  '

Set-TextValue 8 2 'LidoStakedEther'
Set-TextValue 8 3 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextValue 8 4 '3.274.03'
Set-TextValue 8 5 '  -8.88%  '

Set-TextValue 9 2 'USDC'
Set-TextValue 9 3 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue 9 4 '1.00'
Set-TextValue 9 5 '  +0.00%  '

Set-TextValue 10 4 '0.616'
Set-TextValue 10 5 '  -9.85%  '

Set-TextValue 11 4 '57.61'
Set-TextValue 11 5 '  -9.79%  '

Set-TextValue 12 4 '0.130'
Set-TextValue 12 5 '  -11.88%  '

Set-TextValue 13 4 '0.0000253'
Set-TextValue 13 5 '  -10.08%  '

Set-TextValue 14 4 '9.10'
Set-TextValue 14 5 '  -11.16%  '

Set-TextValue 15 4 '3.767.18'
Set-TextValue 15 5 '  -9.89%  '

Set-TextValue 16 4 '0.119'
Set-TextValue 16 5 '  -5.78%  '

Set-TextValue 17 4 '3.258.50'
Set-TextValue 17 5 '  -9.63%  '

Set-TextValue 18 2 'Chainlink'
Set-TextValue 18 3 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 18 4 '17.26'
Set-TextValue 18 5 '  -10.41%  '

Set-TextValue 19 2 'WrappedBTC'
Set-TextValue 19 3 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 19 4 '63.047.15'
Set-TextValue 19 5 '  -7.54%  '

Set-TextValue 20 4 '10.89'
Set-TextValue 20 5 '  -11.08%  '

Set-TextValue 21 4 '0.939'
Set-TextValue 21 5 '  -12.20%  '

Set-TextValue 22 4 '371.90'
Set-TextValue 22 5 '  -8.24%  '

Set-TextValue 23 4 '11.21'
Set-TextValue 23 5 '  -9.75%  '

Set-TextValue 24 4 '79.31'
Set-TextValue 24 5 '  -6.65%  '

Set-TextValue 25 4 '3.63'
Set-TextValue 25 5 '  -13.30%  '

Set-TextValue 26 4 '3.76'
Set-TextValue 26 5 '  -2.53%  '

Set-TextValue 27 5 '  -3.41%  '

Set-TextValue 28 4 '2.63'
Set-TextValue 28 5 '  -9.54%  '

Set-TextValue 29 4 '11.27'
Set-TextValue 29 5 '  -10.15%  '

Set-TextValue 30 4 '8.29'
Set-TextValue 30 5 '  -9.89%  '

Set-TextValue 31 4 '648.19'
Set-TextValue 31 5 '  -10.21%  '

Set-TextValue 32 4 '28.23'
Set-TextValue 32 5 '  -10.69%  '

Set-TextValue 33 4 '6.67'
Set-TextValue 33 5 '  -13.49%  '

Set-TextValue 34 4 '11.12'
Set-TextValue 34 5 '  -8.71%  '

Set-TextValue 35 4 '58.40'
Set-TextValue 35 5 '  -8.14%  '

Set-TextValue 36 5 '  -9.07%  '

Set-TextValue 37 5 '  -0.01%  '

Set-TextValue 38 4 '35.87'
Set-TextValue 38 5 '  -14.48%  '

Set-TextValue 39 4 '0.379'
Set-TextValue 39 5 '  -10.95%  '

Set-TextValue 40 4 '0.996'
Set-TextValue 40 5 '  -0.33%  '

Set-TextValue 41 4 '0.125'
Set-TextValue 41 5 '  -6.06%  '

Set-TextValue 42 4 '28.27'
Set-TextValue 42 5 '  +27.40%  '

Set-TextValue 43 4 '2.878.37'
Set-TextValue 43 5 '  -9.97%  '

Set-TextValue 44 4 '0.0₃0654'
Set-TextValue 44 5 '  -13.07%  '

Set-TextValue 45 4 '2.45'
Set-TextValue 45 5 '  -8.15%  '

Set-TextValue 46 4 '2.64'
Set-TextValue 46 5 '  -20.10%  '

Set-TextValue 47 4 '2.81'
Set-TextValue 47 5 '  +3.91%  '

Set-TextValue 48 4 '2.55'
Set-TextValue 48 5 '  -8.43%  '

Set-TextValue 49 4 '0.0380'
Set-TextValue 49 5 '  -8.38%  '

Set-TextValue 50 2 'Stellar'
Set-TextValue 50 3 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 50 4 '0.123'
Set-TextValue 50 5 '  -6.38%  '

Set-TextValue 51 2 'ApeXProtocol'
Set-TextValue 51 3 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue 51 4 '2.90'
Set-TextValue 51 5 '  -5.81%  '
